$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.019.20"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.869.38"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Formula = "'1.003"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Formula = "'312.21"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Formula = "'1.002"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Formula = "'0.5143"
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("D8").Formula = "'0.3849"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Formula = "'0.08296"
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Formula = "'41.48"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "1.871.03"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Formula = "'7.287"
$ws.Range("D16").Formula = "'1.003"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Formula = "'0.00001098"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Formula = "'90.74"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Formula = "'0.06648"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Formula = "'17.72"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Formula = "'6.029"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "28.057.58"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Formula = "'11.10"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").Formula = "'2.248"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Formula = "'3.386"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.074.20"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Formula = "'2.517"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Formula = "'157.40"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Formula = "'20.54"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Formula = "'124.97"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Formula = "'0.1066"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Formula = "'1.033"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Formula = "'5.824"
$ws.Range("E34").Value = "  +3.93%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Formula = "'3.588"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Formula = "'9.477"
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("D37").Formula = "'0.06528"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Formula = "'0.02416"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Formula = "'0.2197"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Formula = "'0.6553"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Formula = "'1.205"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Formula = "'5.025"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Formula = "'1.211"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Formula = "'11.22"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Formula = "'0.6131"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Formula = "'13.13"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Formula = "'1.278"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Formula = "'3.677"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Formula = "'2.023"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Formula = "'1.215"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Formula = "'120.91"
$ws.Range("E51").Value = "  -0.48%  "
